$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prices")

# Append the new day's row (78) of Solar Prices data, keeping every value
# as plain text (matching the existing inlineStr cells in the sheet).
$row = 78
$values = @("2025-05-18", "37.5", "37", "0.94", "0.258", "0.09", "5,298", "7,931", "7,981", "7.2226")

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 1
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $values[$i]
    $cell.Style = "Normal"
}
